$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers stay the same (row 1 untouched).
# Rows 2-7 correspond to the 6 combinations of
# Sending cluster x Target cluster (ECs/FAPs/sCs x ECs/FAPs)
# for the Gnai2 -> C5ar1 ligand-receptor pair.

$rows = @(
    @{ r=2; A="ECs";  B="Gnai2"; C="C5ar1"; D="ECs";  vals=@(3,1,203.7816646666667,611.344994,0.6667327591988204,0.6667327591988205,2,0.6666666666666666,93.08600066666666,279.258002,0.9956622990949658,0.995662299094966,18969.22017301578,170722.981557142,0.6638406719058277,0.6638406719058279) },
    @{ r=3; A="ECs";  B="Gnai2"; C="C5ar1"; D="FAPs"; vals=@(3,1,203.7816646666667,611.344994,0.6667327591988204,0.6667327591988205,3,1,0.4055383333333333,1.216615,0.004337700905034127,0.004337700905034127,82.64127665281222,743.7714898753101,0.002892087292992623,0.002892087292992624) },
    @{ r=4; A="FAPs"; B="Gnai2"; C="C5ar1"; D="ECs";  vals=@(3,1,63.14058933333333,189.421768,0.2065833519051582,0.2065833519051582,2,0.6666666666666666,93.08600066666666,279.258002,0.9956622990949658,0.995662299094966,5877.504940776393,52897.54446698753,0.2056872551126342,0.2056872551126342) },
    @{ r=5; A="FAPs"; B="Gnai2"; C="C5ar1"; D="FAPs"; vals=@(3,1,63.14058933333333,189.421768,0.2065833519051582,0.2065833519051582,3,1,0.4055383333333333,1.216615,0.004337700905034127,0.004337700905034127,25.60592936392444,230.45336427532,0.000896096792523988,0.0008960967925239884) },
    @{ r=6; A="sCs";  B="Gnai2"; C="C5ar1"; D="ECs";  vals=@(3,1,38.719942,116.159826,0.1266838888960214,0.1266838888960214,2,0.6666666666666666,93.08600066666666,279.258002,0.9956622990949658,0.995662299094966,3604.284546825294,32438.56092142765,0.1261343720765039,0.1261343720765039) },
    @{ r=7; A="sCs";  B="Gnai2"; C="C5ar1"; D="FAPs"; vals=@(3,1,38.719942,116.159826,0.1266838888960214,0.1266838888960214,3,1,0.4055383333333333,1.216615,0.004337700905034127,0.004337700905034127,15.70242074544333,141.32178670899,0.0005495168195175147,0.0005495168195175149) }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    for ($i = 0; $i -lt $row.vals.Count; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $row.vals[$i]
    }
}
